$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This sheet previously contained a duplicated header row (row 2 was
# an exact copy of row 1's column-name headers). The data-cleaning
# fix removes that duplicate row, which shifts the "section title" /
# metric rows up by one, and clears the now-unused trailing row.
# It also clears the bold/border/center formatting that had been
# applied only to row 1, and clears the "Unnamed: 0" label in A1.
# ------------------------------------------------------------------

# Row labels (column A) for rows 2-8 after the shift
$ws.Range("A2").Value = "Fixation based metrics"
$ws.Range("A3").Value = "Revisit count"
$ws.Range("A4").Value = "Fixation count"
$ws.Range("A5").Value = "Dwell time (ms)"
$ws.Range("A6").Value = "Dwell time (%)"
$ws.Range("A7").Value = "Fixation duration (ms)"
$ws.Range("A8").Value = "First fixation duration (ms)"

# Clear A1 ("Unnamed: 0" -> blank)
$ws.Range("A1").Value = ""

# Clear every non-label cell in row 2 (it only carries the section title in A2)
$ws.Range("B2:AR2").ClearContents()

# Numeric metric columns that hold data: G, I, O, U, W, AB, AF, AO
$dataCols = @("G", "I", "O", "U", "W", "AB", "AF", "AO")

# Values per column for rows 3-8 (post-shift), taken from the row
# that used to be one below it before the duplicate header was removed
$values = @{
    "G"  = @(0, 1, 100.06, 0.09, 100.06, 100.06)
    "I"  = @(3, 10, 2903.28, 2.59, 290.33, 100.06)
    "O"  = @(6, 137, 40809.77, 36.4, 297.88, 284.88)
    "U"  = @(0, 1, 600.97, 0.54, 600.97, 600.97)
    "W"  = @(0, 1, 367.07, 0.33, 367.07, 367.07)
    "AB" = @(0, 1, 600.97, 0.54, 600.97, 600.97)
    "AF" = @(9, 89, 27468.28, 24.5, 308.63, 183.49)
    "AO" = @(0, 1, 600.97, 0.54, 600.97, 600.97)
}

foreach ($col in $dataCols) {
    $rowIdx = 0
    for ($r = 3; $r -le 8; $r++) {
        $addr = "$col$r"
        $ws.Range($addr).Value = $values[$col][$rowIdx]
        $rowIdx++
    }
}

# Row 9 no longer holds any data - clear every cell in it
$ws.Range("A9:AR9").ClearContents()

# Remove the bold / centered / thin-border formatting that was only
# ever applied to row 1 (style index 1 in the original file) so every
# cell reverts to the plain default style.
$ws.Cells.ClearFormats()
